# Update Price (D) and Volume(1h) (E) columns with refreshed quotes from the
# symbol-list data pull. Values are stored as literal text (e.g. "322.28",
# "8.44%") in the source workbook, so each target cell is pre-formatted as
# Text before the new value is written - this stops Excel's COM layer from
# auto-coercing numeric-looking strings (or "NN%" strings) into actual
# numbers/percentages, which would change both the stored type and the
# cell's number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","E13","D14","E14","D15","E15","D16","E16","D18","E18","D19","E19","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D38","E38","D39","E39","D40","D41","E41","D42","E42","D43","E43","D44","E44","E45","D46","E46","E47","D49","D50","E50","D51","E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "322.28"
$ws.Range("E2").Value = "8.44%"
$ws.Range("D3").Value = "49.90"
$ws.Range("E3").Value = "19.20%"
$ws.Range("D4").Value = "5.365"
$ws.Range("E4").Value = "7.26%"
$ws.Range("D5").Value = "0.08145"
$ws.Range("E5").Value = "8.33%"
$ws.Range("D6").Value = "4.598"
$ws.Range("E6").Value = "4.89%"
$ws.Range("D7").Value = "1.673"
$ws.Range("E7").Value = "5.68%"
$ws.Range("D8").Value = "1.164"
$ws.Range("E8").Value = "25.66%"
$ws.Range("D9").Value = "0.1353"
$ws.Range("E9").Value = "12.76%"
$ws.Range("D10").Value = "0.1975"
$ws.Range("E10").Value = "7.68%"
$ws.Range("D11").Value = "0.09557"
$ws.Range("E11").Value = "7.19%"
$ws.Range("D12").Value = "0.04583"
$ws.Range("E12").Value = "12.44%"
$ws.Range("E13").Value = "-0.13%"
$ws.Range("D14").Value = "0.001332"
$ws.Range("E14").Value = "3.39%"
$ws.Range("D15").Value = "0.005834"
$ws.Range("E15").Value = "-2.22%"
$ws.Range("D16").Value = "3.386"
$ws.Range("E16").Value = "0.87%"
$ws.Range("D18").Value = "0.3393"
$ws.Range("E18").Value = "2.40%"
$ws.Range("D19").Value = "8.066"
$ws.Range("E19").Value = "-0.54%"
$ws.Range("E20").Value = "1.32%"
$ws.Range("D21").Value = "0.3052"
$ws.Range("E21").Value = "-1.60%"
$ws.Range("D22").Value = "0.04314"
$ws.Range("E22").Value = "4.95%"
$ws.Range("D23").Value = "0.001306"
$ws.Range("E23").Value = "3.20%"
$ws.Range("D24").Value = "0.004308"
$ws.Range("E24").Value = "10.34%"
$ws.Range("D25").Value = "0.0001350"
$ws.Range("E25").Value = "9.69%"
$ws.Range("D26").Value = "0.0003721"
$ws.Range("E26").Value = "-0.07%"
$ws.Range("D38").Value = "0.02763"
$ws.Range("E38").Value = "14.93%"
$ws.Range("D39").Value = "0.05526"
$ws.Range("E39").Value = "5.55%"
$ws.Range("D40").Value = "0.006199"
$ws.Range("D41").Value = "0.007793"
$ws.Range("E41").Value = "-0.42%"
$ws.Range("D42").Value = "0.1448"
$ws.Range("E42").Value = "9.29%"
$ws.Range("D43").Value = "0.007676"
$ws.Range("E43").Value = "3.89%"
$ws.Range("D44").Value = "0.008833"
$ws.Range("E44").Value = "13.01%"
$ws.Range("E45").Value = "18.31%"
$ws.Range("D46").Value = "0.00006763"
$ws.Range("E46").Value = "3.63%"
$ws.Range("E47").Value = "-0.11%"
$ws.Range("D49").Value = "0.003999"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "-0.11%"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "-0.11%"
